$wb = $excel.ActiveWorkbook

# --- Active sheet / selection bookkeeping ---------------------------------
# Previously "Thrusters" (index 5, 0-based activeTab=4) was the active tab
# with selection E17. Now "TTC" (index 2, 0-based activeTab=1) becomes the
# active tab, with selection F13. Activating the TTC sheet and setting its
# selection there will also flip tabSelected appropriately and clear it from
# Thrusters once TTC becomes ActiveSheet.
$wsTTC = $wb.Worksheets.Item("TTC")

# --- TTC sheet data edits ---------------------------------------------------
# Row 2 (TTC-earth): swap face1/face2 (y- <-> z+), swap offset1/offset2
# (0.75 <-> 0.5), simplify the area formula, and add a note in H2.
$wsTTC.Cells.Item(2, 2).Value = "z+"
$wsTTC.Cells.Item(2, 3).Value = "y-"
$wsTTC.Cells.Item(2, 4).Value = 0.5
$wsTTC.Cells.Item(2, 5).Value = 0.75
$wsTTC.Cells.Item(2, 7).Formula = "=1.5*1.5"

$h2 = $wsTTC.Cells.Item(2, 8)
$h2.ClearFormats()
$h2.WrapText = $true
$h2.WrapText = $false
$h2.Value = "Worst case for now. Better:side view for area, rotate behind/front of body for other pointing, but will never present full s/a"

# Row 3 (TTC-nadir): offset1 0.5 -> 0.25
$wsTTC.Cells.Item(3, 4).Value = 0.25

# --- Selection / active tab ------------------------------------------------
$null = $wsTTC.Activate()
$null = $wsTTC.Range("F13").Select()
